$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "Glasgow"
$ws.Range("B12").Value = "Amore"

$ws.Range("A13").Value = "Glasgow"
$ws.Range("B13").Value = "Curlers Rest"

$ws.Range("A14").Value = "Glasgow"
$ws.Range("B14").Value = "Curlers Rest"

$ws.Range("A15").Value = "Glasgow"
$ws.Range("B15").Value = "Curlers Rest"

$ws.Range("A16").Value = "Glasgow"
$ws.Range("B16").Value = "Curlers Rest"

$ws.Range("B16").Select()
